$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 25)
$ws.Range("D2").Value = "4.661320246727603E-13" -as [double]
$ws.Range("E2").Value = "4.661320246727603E-13" -as [double]

# Row 3 (Control 44)
$ws.Range("D3").Value = "3.151212819715415E-42" -as [double]
$ws.Range("E3").Value = "3.151212819715415E-42" -as [double]

# Row 4 (Control 40)
$ws.Range("D4").Value = "4.719402372206117E-26" -as [double]
$ws.Range("E4").Value = "4.719402372206117E-26" -as [double]

# Row 5 (Control 41)
$ws.Range("D5").Value = "7.808695147679701E-09" -as [double]
$ws.Range("E5").Value = "7.808695147679701E-09" -as [double]

# Row 7 (MDD 37)
$ws.Range("D7").Value = "0.9999999999999998" -as [double]
$ws.Range("E7").Value = "2.220446049250313E-16" -as [double]

# Row 8 (MDD 24)
$ws.Range("D8").Value = "5.745772629094356E-31" -as [double]

# Row 9 (MDD 6)
$ws.Range("D9").Value = "6.764328687319019E-25" -as [double]
$ws.Range("E9").Value = 1

# Row 10 (MDD 54)
$ws.Range("D10").Value = "3.441526873966667E-15" -as [double]
$ws.Range("E10").Value = "0.9999999999999966" -as [double]

# Row 11 (MDD 21)
$ws.Range("D11").Value = "8.975262635667823E-12" -as [double]
$ws.Range("E11").Value = "0.9999999999910247" -as [double]
$ws.Range("F11").Value = "24.42014694213867" -as [double]
